$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1084.4445
$ws.Range("I29").Value = 1084.4445
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3253.3335
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2972.3335

$ws.Range("H69").Value = 10665.333
$ws.Range("I69").Value = 11998
$ws.Range("J69").Value = 9999
$ws.Range("K69").Value = 35994
$ws.Range("L69").Value = 29997
$ws.Range("M69").Value = -35120
$ws.Range("N69").Value = -31745

$ws.Range("H72").Value = 10665.333
$ws.Range("I72").Value = 11998
$ws.Range("J72").Value = 9999
$ws.Range("K72").Value = 107982
$ws.Range("L72").Value = 89991
$ws.Range("M72").Value = -103614
$ws.Range("N72").Value = -98727

$ws.Range("H74").Value = 4631
$ws.Range("I74").Value = 4631
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4631
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3695

$ws.Range("H77").Value = 4631
$ws.Range("I77").Value = 4631
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 23155
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -18475

$ws.Range("H132").Value = 6271.636
$ws.Range("I132").Value = 6748.9
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 20246.7
$ws.Range("L132").Value = 4497
$ws.Range("M132").Value = -17716.7
$ws.Range("N132").Value = -9557

$ws.Range("H138").Value = 704.03705
$ws.Range("I138").Value = 697.5833
$ws.Range("J138").Value = 755.6667
$ws.Range("K138").Value = 2092.7499
$ws.Range("L138").Value = 2267.0001
$ws.Range("M138").Value = 3047.2501
$ws.Range("N138").Value = -12547.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 649.5
$ws.Range("I4").Value = 649.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 649.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -533.5

$ws.Range("H6").Value = 1017500
$ws.Range("I6").Value = 1017500
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1017500
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1017327

$ws.Range("H61").Value = 3168.5
$ws.Range("I61").Value = 2947.3635
$ws.Range("J61").Value = 5601
$ws.Range("K61").Value = 2947.3635
$ws.Range("L61").Value = 5601
$ws.Range("M61").Value = -2735.3635
$ws.Range("N61").Value = -6025

$ws.Range("H132").Value = 2070.7856
$ws.Range("I132").Value = 2070.7856
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6212.3568
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3682.3568

$ws.Range("H136").Value = 3168.5
$ws.Range("I136").Value = 2947.3635
$ws.Range("J136").Value = 5601
$ws.Range("K136").Value = 8842.0905
$ws.Range("L136").Value = 16803
$ws.Range("M136").Value = -6292.0905
$ws.Range("N136").Value = -21903

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 368.8
$ws.Range("I22").Value = 368.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 368.8
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -195.8

$ws.Range("H86").Value = 9232.5
$ws.Range("I86").Value = 3374
$ws.Range("J86").Value = 20949.5
$ws.Range("K86").Value = 3374
$ws.Range("L86").Value = 20949.5
$ws.Range("M86").Value = -2251
$ws.Range("N86").Value = -23195.5

$ws.Range("H89").Value = 9232.5
$ws.Range("I89").Value = 3374
$ws.Range("J89").Value = 20949.5
$ws.Range("K89").Value = 16870
$ws.Range("L89").Value = 104747.5
$ws.Range("M89").Value = -11254
$ws.Range("N89").Value = -115979.5

$ws.Range("H134").Value = 5310.05
$ws.Range("I134").Value = 5221.5
$ws.Range("J134").Value = 6107
$ws.Range("K134").Value = 15664.5
$ws.Range("L134").Value = 18321
$ws.Range("M134").Value = -13129.5
$ws.Range("N134").Value = -23391

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3333786
$ws.Range("I22").Value = 531.6
$ws.Range("J22").Value = 5714682
$ws.Range("K22").Value = 531.6
$ws.Range("L22").Value = 5714682
$ws.Range("M22").Value = -181.6
$ws.Range("N22").Value = -5715382

$ws.Range("H37").Value = 24746.75
$ws.Range("I37").Value = 24663
$ws.Range("J37").Value = 24998
$ws.Range("K37").Value = 24663
$ws.Range("L37").Value = 24998
$ws.Range("M37").Value = -24556
$ws.Range("N37").Value = -25212

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 2444
$ws.Range("I82").Value = 2444
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 7332
$ws.Range("L82").ClearContents()
$ws.Range("M82").Value = -6926
$ws.Range("N82").Value = 0

$ws.Range("H85").Value = 2444
$ws.Range("I85").Value = 2444
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 7332
$ws.Range("L85").ClearContents()
$ws.Range("M85").Value = -5928
$ws.Range("N85").Value = 0

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0

$ws.Range("H107").Value = 730.6
$ws.Range("I107").Value = 738.375
$ws.Range("J107").Value = 699.5
$ws.Range("K107").Value = 2215.125
$ws.Range("L107").Value = 2098.5
$ws.Range("M107").Value = -295.125
$ws.Range("N107").Value = -5938.5

$ws.Range("H122").Value = 539.44446
$ws.Range("I122").Value = 494
$ws.Range("J122").Value = 630.3333
$ws.Range("K122").Value = 4446
$ws.Range("L122").Value = 5672.9997
$ws.Range("M122").Value = -1996
$ws.Range("N122").Value = -10572.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 0

$ws.Range("H24").Value = 7499999
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 7499999
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 7499999
$ws.Range("N24").Value = -7500345

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("N26").Value = 0

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("N50").Value = 0

$ws.Range("H126").Value = 2901.375
$ws.Range("I126").Value = 3202
$ws.Range("J126").Value = 1999.5
$ws.Range("K126").Value = 9606
$ws.Range("L126").Value = 5998.5
$ws.Range("M126").Value = -7136
$ws.Range("N126").Value = -10938.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 10325.5
$ws.Range("I56").Value = 14651
$ws.Range("J56").Value = 6000
$ws.Range("K56").Value = 14651
$ws.Range("L56").Value = 6000
$ws.Range("M56").Value = -13960
$ws.Range("N56").Value = -7382

$ws.Range("H61").Value = 2328.875
$ws.Range("I61").Value = 2221.8333
$ws.Range("J61").Value = 2650
$ws.Range("K61").Value = 2221.8333
$ws.Range("L61").Value = 2650
$ws.Range("M61").Value = -2019.8333
$ws.Range("N61").Value = -3054

$ws.Range("H88").Value = 28332.334
$ws.Range("I88").Value = 27499.5
$ws.Range("J88").Value = 29998
$ws.Range("K88").Value = 27499.5
$ws.Range("L88").Value = 29998
$ws.Range("M88").Value = -27071.5
$ws.Range("N88").Value = -30854

$ws.Range("H91").Value = 28332.334
$ws.Range("I91").Value = 27499.5
$ws.Range("J91").Value = 29998
$ws.Range("K91").Value = 27499.5
$ws.Range("L91").Value = 29998
$ws.Range("M91").Value = -26017.5
$ws.Range("N91").Value = -32962

$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -959

$ws.Range("H113").Value = 2328.875
$ws.Range("I113").Value = 2221.8333
$ws.Range("J113").Value = 2650
$ws.Range("K113").Value = 2221.8333
$ws.Range("L113").Value = 2650
$ws.Range("M113").Value = -51.83329999999978
$ws.Range("N113").Value = -6990

$ws.Range("H118").Value = 42500
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 42500
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 42500
$ws.Range("N118").Value = -45814

$ws.Range("H130").Value = 66658.336
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 66658.336
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 66658.336
$ws.Range("N130").Value = -76698.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 225.76471
$ws.Range("I9").Value = 225.76471
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 225.76471
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -85.76471000000001

$ws.Range("H61").Value = 99828.42999999999
$ws.Range("I61").Value = 182666.33
$ws.Range("J61").Value = 37700
$ws.Range("K61").Value = 182666.33
$ws.Range("L61").Value = 37700
$ws.Range("M61").Value = -182374.33
$ws.Range("N61").Value = -38284
